$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# "cambios de las fracciones" - roll the reporting period forward one
# quarter: Q3 2022 (01/07/2022 - 30/09/2022) -> Q4 2022 (01/10/2022 - 31/12/2022),
# and push the validation/update dates out to 10/01/2023.
$ws.Range("B8").Value = 44835   # Fecha de inicio del periodo que se informa -> 01/10/2022
$ws.Range("C8").Value = 44926   # Fecha de término del periodo que se informa -> 31/12/2022
$ws.Range("K8").Value = 44936   # Fecha de validación -> 10/01/2023
$ws.Range("L8").Value = 44936   # Fecha de actualización -> 10/01/2023

# Match the saved view/selection state of the sheet.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 2
$ws.Range("L11").Select()
